{"js": "// Replace the 25 multiplication expressions in the practice-sheet table.\n// Each old expression is unique in the document, so an exact, case-sensitive\n// search-and-replace on the body text is sufficient and safe (the date\n// paragraph above the table is left untouched).\nconst replacements = [\n  [\"929\u00d79=\", \"741\u00d77=\"],\n  [\"497\u00d74=\", \"646\u00d75=\"],\n  [\"601\u00d77=\", \"112\u00d74=\"],\n  [\"709\u00d78=\", \"980\u00d79=\"],\n  [\"753\u00d74=\", \"203\u00d72=\"],\n  [\"901\u00d76=\", \"503\u00d78=\"],\n  [\"608\u00d75=\", \"646\u00d73=\"],\n  [\"994\u00d72=\", \"671\u00d75=\"],\n  [\"912\u00d73=\", \"558\u00d78=\"],\n  [\"148\u00d72=\", \"883\u00d73=\"],\n  [\"208\u00d78=\", \"447\u00d74=\"],\n  [\"651\u00d74=\", \"255\u00d73=\"],\n  [\"424\u00d76=\", \"152\u00d73=\"],\n  [\"597\u00d76=\", \"332\u00d73=\"],\n  [\"626\u00d79=\", \"876\u00d78=\"],\n  [\"457\u00d77=\", \"980\u00d72=\"],\n  [\"529\u00d72=\", \"212\u00d76=\"],\n  [\"397\u00d75=\", \"985\u00d75=\"],\n  [\"558\u00d74=\", \"459\u00d72=\"],\n  [\"441\u00d77=\", \"774\u00d79=\"],\n  [\"874\u00d79=\", \"747\u00d79=\"],\n  [\"582\u00d76=\", \"640\u00d79=\"],\n  [\"962\u00d77=\", \"998\u00d75=\"],\n  [\"800\u00d73=\", \"515\u00d76=\"],\n  [\"525\u00d74=\", \"543\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 multiplication expressions in the practice-sheet table.\n# Each old expression is unique in the document, so a literal (non-wildcard)\n# Find/Replace-All on the whole document content is sufficient and safe\n# (the date paragraph above the table is left untouched since it does not\n# match any of the search strings).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"929\u00d79=\", \"741\u00d77=\"),\n    @(\"497\u00d74=\", \"646\u00d75=\"),\n    @(\"601\u00d77=\", \"112\u00d74=\"),\n    @(\"709\u00d78=\", \"980\u00d79=\"),\n    @(\"753\u00d74=\", \"203\u00d72=\"),\n    @(\"901\u00d76=\", \"503\u00d78=\"),\n    @(\"608\u00d75=\", \"646\u00d73=\"),\n    @(\"994\u00d72=\", \"671\u00d75=\"),\n    @(\"912\u00d73=\", \"558\u00d78=\"),\n    @(\"148\u00d72=\", \"883\u00d73=\"),\n    @(\"208\u00d78=\", \"447\u00d74=\"),\n    @(\"651\u00d74=\", \"255\u00d73=\"),\n    @(\"424\u00d76=\", \"152\u00d73=\"),\n    @(\"597\u00d76=\", \"332\u00d73=\"),\n    @(\"626\u00d79=\", \"876\u00d78=\"),\n    @(\"457\u00d77=\", \"980\u00d72=\"),\n    @(\"529\u00d72=\", \"212\u00d76=\"),\n    @(\"397\u00d75=\", \"985\u00d75=\"),\n    @(\"558\u00d74=\", \"459\u00d72=\"),\n    @(\"441\u00d77=\", \"774\u00d79=\"),\n    @(\"874\u00d79=\", \"747\u00d79=\"),\n    @(\"582\u00d76=\", \"640\u00d79=\"),\n    @(\"962\u00d77=\", \"998\u00d75=\"),\n    @(\"800\u00d73=\", \"515\u00d76=\"),\n    @(\"525\u00d74=\", \"543\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute([Type]::Missing, $true, $false, $false, $false, $false, $true, 1, $false, [Type]::Missing, 2) | Out-Null\n}\n"}
